# Add administrative geolocation levels to the SpecimenInfo header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SpecimenInfo")

# Full set of header labels for row 1, in column order (A..U).
# This inserts collection_country, geo_admin1, geo_admin2, geo_admin3
# right before lat_lon, and drops the old geo_loc_name column,
# shifting every subsequent header three columns to the right.
$headers = @(
    "specimen_id",
    "plate_name",
    "plate_row",
    "plate_col",
    "samp_taxon_id",
    "individual_id",
    "host_taxon_id",
    "alternate_identifiers",
    "parasite_density",
    "collection_date",
    "collection_country",
    "geo_admin1",
    "geo_admin2",
    "geo_admin3",
    "lat_lon",
    "collector",
    "samp_store_loc",
    "samp_collect_device",
    "project_name",
    "accession",
    "sample_comments"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
